# Fill in the W4 Salaries and Tasks template with the team's actual data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header info -----------------------------------------------------
$ws.Range("B3").Value = 44112          # Date -> 10/8/2020
$ws.Range("B4").Value = "Limette"      # Team Name
$ws.Range("B5").Value = 4              # Total Number of Team Members

# --- Team member names + salaries ------------------------------------
$ws.Range("A8").Value  = "Lukas Hasler"
$ws.Range("B8").Value  = 100

$ws.Range("A9").Value  = "Pascal Strebel"
$ws.Range("B9").Value  = 100

$ws.Range("A10").Value = "Cedric Weibel"
$ws.Range("B10").Value = 100

$ws.Range("A11").Value = "Robin Schmidiger"
$ws.Range("B11").Value = 100

# 5th member row is left blank (only 4 team members)
$ws.Range("A12").Value = $null

# --- Row 18 shrinks now that the task instructions wrap less ---------
$ws.Rows.Item(18).RowHeight = 39

# --- Tasks completed this week / to complete next week ---------------
$ws.Range("A19").Value = "Brainstorming"
$ws.Range("A20").Value = "Pesonas"
$ws.Range("A21").Value = "Presentation"
$ws.Range("B19").Value = "Prototypes"

# --- Final selection matches the saved workbook state -----------------
$ws.Range("B20").Select() | Out-Null
